$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value while keeping it stored as text (these tables
# store ratios/counts as text strings, e.g. "0.53", not as numbers), and
# without leaving the cell with a non-default (Text) number format/style.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 3 (Female, Age 60 - 69)
Set-TextValue "C3" "0.53"
Set-TextValue "E3" "0.5"
Set-TextValue "G3" "0.52"
Set-TextValue "I3" "0.5"

# Row 4 (Female, Age 70 - 79)
Set-TextValue "C4" "0.31"
Set-TextValue "E4" "0.32"
Set-TextValue "G4" "0.31"
Set-TextValue "I4" "0.31"

# Row 5 (Female, Age 80 - 89)
Set-TextValue "C5" "0.13"
Set-TextValue "E5" "0.15"
Set-TextValue "G5" "0.14"
Set-TextValue "I5" "0.15"

# Row 6 (Female, Age 90 plus)
Set-TextValue "C6" "0.03"
Set-TextValue "E6" "0.03"
Set-TextValue "G6" "0.04"
Set-TextValue "I6" "0.04"

# Row 8 (Female, Less than Primary)
Set-TextValue "C8" "0.61"
Set-TextValue "E8" "0.26"
Set-TextValue "I8" "0.27"

# Row 9 (Female, Primary)
Set-TextValue "C9" "0.29"
Set-TextValue "G9" "0.18"

# Row 10 (Female, Secondary)
Set-TextValue "C10" "0.07"
Set-TextValue "D10" "0.21"
Set-TextValue "E10" "0.36"
Set-TextValue "F10" "0.42"
Set-TextValue "H10" "0.28"
Set-TextValue "J10" "0.45"

# Row 11 (Female, University)
Set-TextValue "G11" "0.04"

# Row 13 (Female, Household Size)
Set-TextValue "C13" "3.56"
Set-TextValue "G13" "3.55"
Set-TextValue "I13" "3.01"

# Row 14 (Female, Lives Alone)
Set-TextValue "C14" "0.13"
Set-TextValue "E14" "0.27"
Set-TextValue "G14" "0.12"

# Row 15 (Female, Lives with Child)
Set-TextValue "E15" "0.29"

# Row 16 (Female, Married/Cohabiting)
Set-TextValue "C16" "0.46"
Set-TextValue "E16" "0.42"
Set-TextValue "G16" "0.39"
Set-TextValue "I16" "0.46"

# Row 19 (Male, Age 60 - 69)
Set-TextValue "E19" "0.54"
Set-TextValue "G19" "0.54"
Set-TextValue "I19" "0.52"

# Row 20 (Male, Age 70 - 79)
Set-TextValue "C20" "0.31"
Set-TextValue "E20" "0.32"
Set-TextValue "G20" "0.3"
Set-TextValue "I20" "0.32"

# Row 21 (Male, Age 80 - 89)
Set-TextValue "C21" "0.12"
Set-TextValue "E21" "0.13"
Set-TextValue "G21" "0.12"

# Row 24 (Male, Less than Primary)
Set-TextValue "I24" "0.21"

# Row 26 (Male, Secondary)
Set-TextValue "D26" "0.21"
Set-TextValue "F26" "0.42"
Set-TextValue "H26" "0.32"
Set-TextValue "I26" "0.23"
Set-TextValue "J26" "0.45"

# Row 27 (Male, University)
Set-TextValue "I27" "0.08"

# Row 30 (Male, Lives Alone)
Set-TextValue "E30" "0.18"
